$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 1, shifting existing data down.
$ws.Rows("1:1").Insert()

# Populate the new header row with the customer's priority number.
$ws.Range("A1").Value = "Priority"
$ws.Range("B1").Value = 1

# Match the selection left by the editor.
$ws.Range("B1").Select()
